$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A7").Value = 1
$ws.Range("A17").Value = 0
$ws.Range("A18").Value = 0
$ws.Range("A31").Value = 0
$ws.Range("A34").Value = 0
$ws.Range("A64").Value = 0
$ws.Range("A70").Value = 0
$ws.Range("A109").Value = 0
$ws.Range("A110").Value = 1
$ws.Range("A121").Value = 1
$ws.Range("A140").Value = 0
$ws.Range("A149").Value = 1
$ws.Range("A152").Value = 1
$ws.Range("A164").Value = 0
$ws.Range("A173").Value = 0
$ws.Range("A174").Value = 0
$ws.Range("A175").Value = 0
$ws.Range("A176").Value = 0
$ws.Range("A443").Value = 0
$ws.Range("A655").Value = 0
$ws.Range("A944").Value = 1
$ws.Range("A974").Value = 1
$ws.Range("A1046").Value = 1
$ws.Range("A1055").Value = 0
$ws.Range("A1084").Value = 1
$ws.Range("A1177").Value = 1
$ws.Range("A1202").Value = 1
$ws.Range("A1213").Value = 0
$ws.Range("A1215").Value = 0
$ws.Range("A1219").Value = 0
$ws.Range("A1220").Value = 0
